# Auto-generated Excel COM-interop script to apply odds updates
# as described by the diff for Jogos_da_Semana_FlashScore_2025-02-17.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Cells.Item(2, 14).Value = 9   # N2: 8.5 -> 9
$ws.Cells.Item(2, 23).Value = 1.91   # W2: 1.95 -> 1.91
$ws.Cells.Item(2, 24).Value = 1.91   # X2: 1.8 -> 1.91
$ws.Cells.Item(2, 31).Value = 9   # AE2: 8.5 -> 9
$ws.Cells.Item(2, 36).Value = 12   # AJ2: 11 -> 12

# Row 3
$ws.Cells.Item(3, 10).Value = 1.57   # J3: 1.54 -> 1.57
$ws.Cells.Item(3, 14).Value = 23   # N3: 26 -> 23

# Row 4
$ws.Cells.Item(4, 10).Value = 2.63   # J4: 2.62 -> 2.63
$ws.Cells.Item(4, 11).Value = 1.91   # K4: 1.87 -> 1.91
$ws.Cells.Item(4, 13).Value = 1.13   # M4: 1.1 -> 1.13
$ws.Cells.Item(4, 15).Value = 1.57   # O4: 1.54 -> 1.57
$ws.Cells.Item(4, 18).Value = 1.44   # R4: 1.41 -> 1.44
$ws.Cells.Item(4, 20).Value = 1.13   # T4: 1.1 -> 1.13

# Row 5
$ws.Cells.Item(5, 7).Value = 2.88   # G5: 2.8 -> 2.88
$ws.Cells.Item(5, 9).Value = 2.9   # I5: 3 -> 2.9
$ws.Cells.Item(5, 11).Value = 1.83   # K5: 1.8 -> 1.83
$ws.Cells.Item(5, 13).Value = 1.13   # M5: 1.1 -> 1.13
$ws.Cells.Item(5, 15).Value = 1.57   # O5: 1.54 -> 1.57
$ws.Cells.Item(5, 17).Value = 2.88   # Q5: 2.87 -> 2.88
$ws.Cells.Item(5, 18).Value = 1.4   # R5: 1.37 -> 1.4
$ws.Cells.Item(5, 20).Value = 1.13   # T5: 1.1 -> 1.13
$ws.Cells.Item(5, 21).Value = 1.62   # U5: 1.67 -> 1.62
$ws.Cells.Item(5, 22).Value = 2.2   # V5: 2.1 -> 2.2
$ws.Cells.Item(5, 37).Value = 12   # AK5: 13 -> 12
$ws.Cells.Item(5, 39).Value = 29   # AM5: 34 -> 29

# Row 6
$ws.Cells.Item(6, 11).Value = 1.8   # K6: 1.77 -> 1.8
$ws.Cells.Item(6, 13).Value = 1.17   # M6: 1.13 -> 1.17
$ws.Cells.Item(6, 15).Value = 1.73   # O6: 1.69 -> 1.73
$ws.Cells.Item(6, 18).Value = 1.3   # R6: 1.27 -> 1.3
$ws.Cells.Item(6, 20).Value = 1.08   # T6: 1.05 -> 1.08

# Row 7
$ws.Cells.Item(7, 11).Value = 1.8   # K7: 1.77 -> 1.8
$ws.Cells.Item(7, 13).Value = 1.14   # M7: 1.11 -> 1.14
$ws.Cells.Item(7, 15).Value = 1.73   # O7: 1.69 -> 1.73
$ws.Cells.Item(7, 18).Value = 1.33   # R7: 1.3 -> 1.33
$ws.Cells.Item(7, 20).Value = 1.08   # T7: 1.05 -> 1.08

# Row 8
$ws.Cells.Item(8, 13).Value = 1.08   # M8: 1.05 -> 1.08
$ws.Cells.Item(8, 15).Value = 1.44   # O8: 1.41 -> 1.44
$ws.Cells.Item(8, 16).Value = 2.63   # P8: 2.62 -> 2.63
$ws.Cells.Item(8, 18).Value = 1.53   # R8: 1.5 -> 1.53
$ws.Cells.Item(8, 20).Value = 1.17   # T8: 1.13 -> 1.17

# Row 9
$ws.Cells.Item(9, 7).Value = 2.15   # G9: 2.2 -> 2.15
$ws.Cells.Item(9, 9).Value = 3.9   # I9: 3.8 -> 3.9
$ws.Cells.Item(9, 11).Value = 1.83   # K9: 1.8 -> 1.83
$ws.Cells.Item(9, 13).Value = 1.14   # M9: 1.1 -> 1.14
$ws.Cells.Item(9, 14).Value = 5.5   # N9: 6 -> 5.5
$ws.Cells.Item(9, 15).Value = 1.57   # O9: 1.54 -> 1.57
$ws.Cells.Item(9, 17).Value = 2.88   # Q9: 2.87 -> 2.88
$ws.Cells.Item(9, 18).Value = 1.4   # R9: 1.37 -> 1.4
$ws.Cells.Item(9, 20).Value = 1.13   # T9: 1.1 -> 1.13
$ws.Cells.Item(9, 21).Value = 1.67   # U9: 1.62 -> 1.67
$ws.Cells.Item(9, 22).Value = 2.1   # V9: 2.2 -> 2.1
$ws.Cells.Item(9, 23).Value = 2.38   # W9: 2.25 -> 2.38
$ws.Cells.Item(9, 24).Value = 1.53   # X9: 1.57 -> 1.53
$ws.Cells.Item(9, 28).Value = 19   # AB9: 21 -> 19

# Row 10
$ws.Cells.Item(10, 7).Value = 1.85   # G10: 1.9 -> 1.85
$ws.Cells.Item(10, 8).Value = 3   # H10: 2.9 -> 3
$ws.Cells.Item(10, 10).Value = 2.63   # J10: 2.75 -> 2.63
$ws.Cells.Item(10, 11).Value = 1.83   # K10: 1.8 -> 1.83
$ws.Cells.Item(10, 14).Value = 5   # N10: 4.75 -> 5
$ws.Cells.Item(10, 15).Value = 1.67   # O10: 1.73 -> 1.67
$ws.Cells.Item(10, 16).Value = 2.1   # P10: 2 -> 2.1
$ws.Cells.Item(10, 17).Value = 3.4   # Q10: 3.5 -> 3.4
$ws.Cells.Item(10, 18).Value = 1.33   # R10: 1.3 -> 1.33
$ws.Cells.Item(10, 19).Value = 7   # S10: 8 -> 7
$ws.Cells.Item(10, 20).Value = 1.1   # T10: 1.08 -> 1.1
$ws.Cells.Item(10, 21).Value = 1.73   # U10: 1.75 -> 1.73
$ws.Cells.Item(10, 22).Value = 2   # V10: 2.05 -> 2
$ws.Cells.Item(10, 25).Value = 4.5   # Y10: 4.33 -> 4.5
$ws.Cells.Item(10, 31).Value = 5   # AE10: 4.75 -> 5
$ws.Cells.Item(10, 33).Value = 26   # AG10: 29 -> 26
$ws.Cells.Item(10, 34).Value = 126   # AH10: 151 -> 126
$ws.Cells.Item(10, 41).Value = 67   # AO10: 81 -> 67

# Row 16
$ws.Cells.Item(16, 15).Value = 1.3   # O16: 1.29 -> 1.3
$ws.Cells.Item(16, 16).Value = 3.4   # P16: 3.5 -> 3.4
$ws.Cells.Item(16, 17).Value = 2   # Q16: 1.98 -> 2
$ws.Cells.Item(16, 18).Value = 1.85   # R16: 1.88 -> 1.85

# Row 17
$ws.Cells.Item(17, 11).Value = 1.91   # K17: 1.87 -> 1.91

# Row 18
$ws.Cells.Item(18, 11).Value = 1.91   # K18: 1.87 -> 1.91

# Row 19
$ws.Cells.Item(19, 12).Value = 2.88   # L19: 2.87 -> 2.88
$ws.Cells.Item(19, 17).Value = 1.8   # Q19: 1.75 -> 1.8
$ws.Cells.Item(19, 18).Value = 2   # R19: 2.05 -> 2

# Row 20
$ws.Cells.Item(20, 7).Value = 1.6   # G20: 1.55 -> 1.6
$ws.Cells.Item(20, 8).Value = 4   # H20: 4.1 -> 4
$ws.Cells.Item(20, 9).Value = 5   # I20: 5.25 -> 5
$ws.Cells.Item(20, 11).Value = 2.38   # K20: 2.4 -> 2.38
$ws.Cells.Item(20, 12).Value = 5   # L20: 5.5 -> 5
$ws.Cells.Item(20, 13).Value = 1.04   # M20: 1.03 -> 1.04
$ws.Cells.Item(20, 14).Value = 13   # N20: 15 -> 13
$ws.Cells.Item(20, 21).Value = 1.33   # U20: 1.3 -> 1.33
$ws.Cells.Item(20, 22).Value = 3.25   # V20: 3.4 -> 3.25

# Row 29
$ws.Cells.Item(29, 11).Value = 2.63   # K29: 2.62 -> 2.63
$ws.Cells.Item(29, 23).Value = 1.44   # W29: 1.41 -> 1.44
$ws.Cells.Item(29, 24).Value = 2.63   # X29: 2.62 -> 2.63

# Row 30
$ws.Cells.Item(30, 23).Value = 1.62   # W30: 1.58 -> 1.62

# Row 31
$ws.Cells.Item(31, 23).Value = 1.83   # W31: 1.8 -> 1.83
$ws.Cells.Item(31, 24).Value = 1.83   # X31: 1.8 -> 1.83
$ws.Cells.Item(31, 45).Value = 2.29   # AS31: 2.28 -> 2.29

# Row 32
$ws.Cells.Item(32, 14).Value = 13   # N32: 12 -> 13
$ws.Cells.Item(32, 15).Value = 1.25   # O32: 1.22 -> 1.25
$ws.Cells.Item(32, 16).Value = 3.75   # P32: 4 -> 3.75
$ws.Cells.Item(32, 17).Value = 1.83   # Q32: 1.8 -> 1.83
$ws.Cells.Item(32, 18).Value = 2.03   # R32: 2 -> 2.03
$ws.Cells.Item(32, 19).Value = 3   # S32: 2.75 -> 3
$ws.Cells.Item(32, 20).Value = 1.36   # T32: 1.4 -> 1.36
$ws.Cells.Item(32, 23).Value = 1.67   # W32: 1.63 -> 1.67

# Row 33
$ws.Cells.Item(33, 23).Value = 1.8   # W33: 1.77 -> 1.8
$ws.Cells.Item(33, 24).Value = 1.91   # X33: 1.87 -> 1.91

# Row 34
$ws.Cells.Item(34, 18).Value = 1.5   # R34: 1.47 -> 1.5
$ws.Cells.Item(34, 24).Value = 1.7   # X34: 1.67 -> 1.7

# Row 35
$ws.Cells.Item(35, 18).Value = 1.67   # R35: 1.63 -> 1.67
$ws.Cells.Item(35, 23).Value = 1.91   # W35: 1.87 -> 1.91
$ws.Cells.Item(35, 24).Value = 1.8   # X35: 1.77 -> 1.8

# Row 37
$ws.Cells.Item(37, 7).Value = 4.45   # G37: 4.1 -> 4.45
$ws.Cells.Item(37, 8).Value = 3.05   # H37: 3.1 -> 3.05
$ws.Cells.Item(37, 9).Value = 1.83   # I37: 1.87 -> 1.83
$ws.Cells.Item(37, 10).Value = 5.1   # J37: 4.6 -> 5.1
$ws.Cells.Item(37, 11).Value = 1.95   # K37: 2.02 -> 1.95
$ws.Cells.Item(37, 12).Value = 2.47   # L37: 2.5 -> 2.47
$ws.Cells.Item(37, 13).Value = 1.11   # M37: 1.09 -> 1.11
$ws.Cells.Item(37, 14).Value = 5.6   # N37: 6.2 -> 5.6
$ws.Cells.Item(37, 15).Value = 1.5   # O37: 1.4 -> 1.5
$ws.Cells.Item(37, 16).Value = 2.42   # P37: 2.72 -> 2.42
$ws.Cells.Item(37, 17).Value = 2.42   # Q37: 2.18 -> 2.42
$ws.Cells.Item(37, 18).Value = 1.5   # R37: 1.6 -> 1.5
$ws.Cells.Item(37, 19).Value = 4.35   # S37: 3.8 -> 4.35
$ws.Cells.Item(37, 20).Value = 1.18   # T37: 1.23 -> 1.18
$ws.Cells.Item(37, 21).Value = 1.53   # U37: 1.45 -> 1.53
$ws.Cells.Item(37, 22).Value = 2.35   # V37: 2.55 -> 2.35
$ws.Cells.Item(37, 23).Value = 2.18   # W37: 1.98 -> 2.18
$ws.Cells.Item(37, 24).Value = 1.62   # X37: 1.75 -> 1.62
$ws.Cells.Item(37, 25).Value = 9.25   # Y37: 10.25 -> 9.25
$ws.Cells.Item(37, 27).Value = 16   # AA37: 14 -> 16
$ws.Cells.Item(37, 28).Value = 90   # AB37: 70 -> 90
$ws.Cells.Item(37, 29).Value = 60   # AC37: 45 -> 60
$ws.Cells.Item(37, 30).Value = 75   # AD37: 55 -> 75
$ws.Cells.Item(37, 31).Value = 5.6   # AE37: 6.2 -> 5.6
$ws.Cells.Item(37, 33).Value = 20   # AG37: 16.5 -> 20
$ws.Cells.Item(37, 34).Value = 150   # AH37: 90 -> 150
$ws.Cells.Item(37, 36).Value = 5.2   # AJ37: 5.9 -> 5.2
$ws.Cells.Item(37, 37).Value = 7.4   # AK37: 8 -> 7.4
$ws.Cells.Item(37, 38).Value = 8.75   # AL37: 8.5 -> 8.75
$ws.Cells.Item(37, 39).Value = 15   # AM37: 16 -> 15
$ws.Cells.Item(37, 40).Value = 18   # AN37: 17 -> 18
$ws.Cells.Item(37, 41).Value = 40   # AO37: 32 -> 40

# Row 38
$ws.Cells.Item(38, 7).Value = 2.22   # G38: 2.12 -> 2.22
$ws.Cells.Item(38, 9).Value = 3.25   # I38: 3.45 -> 3.25
$ws.Cells.Item(38, 10).Value = 2.8   # J38: 2.7 -> 2.8
$ws.Cells.Item(38, 12).Value = 3.7   # L38: 3.85 -> 3.7
$ws.Cells.Item(38, 15).Value = 1.36   # O38: 1.35 -> 1.36
$ws.Cells.Item(38, 16).Value = 2.9   # P38: 2.95 -> 2.9
$ws.Cells.Item(38, 17).Value = 2.05   # Q38: 2.02 -> 2.05
$ws.Cells.Item(38, 19).Value = 3.5   # S38: 3.45 -> 3.5
$ws.Cells.Item(38, 20).Value = 1.26   # T38: 1.27 -> 1.26
$ws.Cells.Item(38, 21).Value = 1.42   # U38: 1.4 -> 1.42
$ws.Cells.Item(38, 22).Value = 2.67   # V38: 2.7 -> 2.67
$ws.Cells.Item(38, 25).Value = 7   # Y38: 6.8 -> 7
$ws.Cells.Item(38, 26).Value = 10.25   # Z38: 9.75 -> 10.25
$ws.Cells.Item(38, 27).Value = 9   # AA38: 8.75 -> 9
$ws.Cells.Item(38, 28).Value = 22   # AB38: 20 -> 22
$ws.Cells.Item(38, 29).Value = 19   # AC38: 18 -> 19
$ws.Cells.Item(38, 32).Value = 6   # AF38: 6.1 -> 6
$ws.Cells.Item(38, 36).Value = 9.25   # AJ38: 9.75 -> 9.25
$ws.Cells.Item(38, 37).Value = 16.5   # AK38: 18.5 -> 16.5
$ws.Cells.Item(38, 38).Value = 11.25   # AL38: 11.5 -> 11.25
$ws.Cells.Item(38, 39).Value = 45   # AM38: 50 -> 45
$ws.Cells.Item(38, 40).Value = 29   # AN38: 30 -> 29

# Row 39
$ws.Cells.Item(39, 23).Value = 1.8   # W39: 1.77 -> 1.8
$ws.Cells.Item(39, 24).Value = 1.91   # X39: 1.87 -> 1.91

# Row 41
$ws.Cells.Item(41, 13).Value = 1.02   # M41: 1.01 -> 1.02
$ws.Cells.Item(41, 15).Value = 1.14   # O41: 1.11 -> 1.14
$ws.Cells.Item(41, 20).Value = 1.67   # T41: 1.63 -> 1.67

# Row 42
$ws.Cells.Item(42, 13).Value = 1.04   # M42: 1.03 -> 1.04
$ws.Cells.Item(42, 15).Value = 1.22   # O42: 1.19 -> 1.22
$ws.Cells.Item(42, 20).Value = 1.4   # T42: 1.37 -> 1.4

# Row 43
$ws.Cells.Item(43, 13).Value = 1.05   # M43: 1.03 -> 1.05
$ws.Cells.Item(43, 15).Value = 1.3   # O43: 1.27 -> 1.3
$ws.Cells.Item(43, 20).Value = 1.29   # T43: 1.25 -> 1.29

# Row 44
$ws.Cells.Item(44, 13).Value = 1.04   # M44: 1.03 -> 1.04
$ws.Cells.Item(44, 15).Value = 1.25   # O44: 1.22 -> 1.25
$ws.Cells.Item(44, 20).Value = 1.36   # T44: 1.33 -> 1.36

# Row 45
$ws.Cells.Item(45, 13).Value = 1.06   # M45: 1.04 -> 1.06
$ws.Cells.Item(45, 15).Value = 1.3   # O45: 1.27 -> 1.3
$ws.Cells.Item(45, 20).Value = 1.29   # T45: 1.25 -> 1.29
$ws.Cells.Item(45, 23).Value = 1.77   # W45: 1.8 -> 1.77
$ws.Cells.Item(45, 24).Value = 1.87   # X45: 1.91 -> 1.87

# Row 46
$ws.Cells.Item(46, 7).Value = 2.6   # G46: 2.62 -> 2.6
$ws.Cells.Item(46, 8).Value = 3.1   # H46: 3 -> 3.1
$ws.Cells.Item(46, 9).Value = 2.62   # I46: 2.65 -> 2.62
$ws.Cells.Item(46, 10).Value = 3.15   # J46: 3.25 -> 3.15
$ws.Cells.Item(46, 11).Value = 2.02   # K46: 1.98 -> 2.02
$ws.Cells.Item(46, 15).Value = 1.36   # O46: 1.4 -> 1.36
$ws.Cells.Item(46, 16).Value = 2.67   # P46: 2.52 -> 2.67
$ws.Cells.Item(46, 17).Value = 2.05   # Q46: 2.15 -> 2.05
$ws.Cells.Item(46, 18).Value = 1.62   # R46: 1.55 -> 1.62
$ws.Cells.Item(46, 19).Value = 3.35   # S46: 3.6 -> 3.35
$ws.Cells.Item(46, 20).Value = 1.23   # T46: 1.2 -> 1.23
$ws.Cells.Item(46, 21).Value = 1.42   # U46: 1.45 -> 1.42
$ws.Cells.Item(46, 22).Value = 2.47   # V46: 2.37 -> 2.47
$ws.Cells.Item(46, 23).Value = 1.78   # W46: 1.85 -> 1.78
$ws.Cells.Item(46, 24).Value = 1.82   # X46: 1.75 -> 1.82
$ws.Cells.Item(46, 25).Value = 7.7   # Y46: 7.1 -> 7.7
$ws.Cells.Item(46, 26).Value = 12.5   # Z46: 12 -> 12.5
$ws.Cells.Item(46, 27).Value = 9.75   # AA46: 10.25 -> 9.75
$ws.Cells.Item(46, 28).Value = 29   # AB46: 30 -> 29
$ws.Cells.Item(46, 29).Value = 23   # AC46: 25 -> 23
$ws.Cells.Item(46, 30).Value = 35   # AD46: 40 -> 35
$ws.Cells.Item(46, 31).Value = 8.25   # AE46: 7.5 -> 8.25
$ws.Cells.Item(46, 32).Value = 6   # AF46: 5.9 -> 6
$ws.Cells.Item(46, 33).Value = 14.5   # AG46: 15.5 -> 14.5
$ws.Cells.Item(46, 34).Value = 75   # AH46: 90 -> 75
$ws.Cells.Item(46, 35).Value = 700   # AI46: 800 -> 700
$ws.Cells.Item(46, 36).Value = 7.5   # AJ46: 7.3 -> 7.5
$ws.Cells.Item(46, 40).Value = 24   # AN46: 25 -> 24
